# Apply "Saving things for later" content edits to the card-data table
# on the "Sheet1" worksheet (the sheet that holds the game card rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Invasive species card (row 2) - rewritten instructions text
$ws.Range("C2").Value = "A chosen player loses 1 Flora or Fauna each turn. 2 facilities cards needed to clear it. Card stays in front of selected player until cleared."

# Bushfire card (row 3) - add explicit dice-roll wording
$ws.Range("C3").Value = "The player you select rolls 1 dice. Roll 1-5: they must discard all Flora cards. Roll 6: the plague affects you instead."

# Feral Animals card (row 4) - add explicit dice-roll wording
$ws.Range("C4").Value = "The player you select rolls 1 dice. Roll 1-5: they must discard all Fauna cards. Roll 6: the plague affects you instead."

# Bushfire Disaster card (row 7)
$ws.Range("B7").Value = " Emergency"
$ws.Range("C7").Value = "State Emergency! All players (including yourself) lose all Flora cards."

# Lightning Storm card (row 8)
$ws.Range("C8").Value = "A selected player rolls 1 dice. Roll 1-3: they give that number of Fauna Cards to you. Roll 4-6: they keep their hand."

# Duplicate "Invasive species" row (row 9) mirrors row 2's updated text
$ws.Range("C9").Value = "A chosen player loses 1 Flora or Fauna each turn. 2 facilities cards needed to clear it. Card stays in front of selected player until cleared."

# Duplicate "Feral Animals" row (row 10) mirrors row 4's updated text
$ws.Range("C10").Value = "The player you select rolls 1 dice. Roll 1-5: they must discard all Fauna cards. Roll 6: the plague affects you instead."
